# Refresh crypto price/volume snapshot (GitHub Actions scrape update).
# Numeric-looking price strings are quote-prefixed so Excel stores them
# as literal text (matching the source data) instead of coercing them
# to numbers, which would drop trailing zeros / introduce float noise.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.339.44'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '3.686.94'
$ws.Range('E3').Value = '  -3.15%  '
$ws.Range('D5').Value = '''683.85'
$ws.Range('E5').Value = '  -3.54%  '
$ws.Range('D6').Value = '''162.48'
$ws.Range('E6').Value = '  -4.75%  '
$ws.Range('D7').Value = '3.684.28'
$ws.Range('E7').Value = '  -3.23%  '
$ws.Range('E9').Value = '  -4.29%  '
$ws.Range('E10').Value = '  -7.69%  '
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('D12').Value = '''0.450'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  -6.33%  '
$ws.Range('D14').Value = '''33.62'
$ws.Range('D15').Value = '4.309.31'
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('D16').Value = '3.692.64'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').Value = '69.400.27'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').Value = '''16.33'
$ws.Range('E19').Value = '  -6.01%  '
$ws.Range('E20').Value = '  -6.81%  '
$ws.Range('D21').Value = '''482.29'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').Value = '''9.81'
$ws.Range('E22').Value = '  -7.90%  '
$ws.Range('D23').Value = '''0.668'
$ws.Range('E23').Value = '  -8.44%  '
$ws.Range('D24').Value = '''80.00'
$ws.Range('E24').Value = '  -5.24%  '
$ws.Range('D25').Value = '3.833.14'
$ws.Range('E25').Value = '  -3.15%  '
$ws.Range('E26').Value = '  -11.08%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''11.49'
$ws.Range('E28').Value = '  -5.07%  '
$ws.Range('E29').Value = '  -7.88%  '
$ws.Range('E30').Value = '  -10.54%  '
$ws.Range('E31').Value = '  -10.64%  '
$ws.Range('E32').Value = '  -5.78%  '
$ws.Range('E33').Value = '  -6.88%  '
$ws.Range('D34').Value = '''27.09'
$ws.Range('E34').Value = '  -6.88%  '
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('D37').Value = '3.650.76'
$ws.Range('E37').Value = '  -3.36%  '
$ws.Range('E38').Value = '  -6.03%  '
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('D40').Value = '''0.0946'
$ws.Range('E40').Value = '  -6.88%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '''2.17'
$ws.Range('E42').Value = '  -5.98%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -7.58%  '
$ws.Range('D45').Value = '''158.11'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').Value = '''48.16'
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('E47').Value = '  -12.78%  '
$ws.Range('E48').Value = '  -13.60%  '
$ws.Range('E49').Value = '  -3.55%  '
$ws.Range('D50').Value = '''390.19'
$ws.Range('E50').Value = '  -8.31%  '
$ws.Range('E51').Value = '  -5.89%  '
